$d = $word.ActiveDocument

$replacements = @(
    @{old = "84×43=3612"; new = "58×87=5046"},
    @{old = "75×77=5775"; new = "40×31=1240"},
    @{old = "69×86=5934"; new = "50×21=1050"},
    @{old = "74×34=2516"; new = "29×21=609"},
    @{old = "40×27=1080"; new = "93×33=3069"},
    @{old = "64×41=2624"; new = "31×60=1860"},
    @{old = "37×40=1480"; new = "91×96=8736"},
    @{old = "16×28=448"; new = "40×58=2320"},
    @{old = "41×85=3485"; new = "82×24=1968"},
    @{old = "48×63=3024"; new = "41×86=3526"},
    @{old = "23×19=437"; new = "86×23=1978"},
    @{old = "79×86=6794"; new = "89×40=3560"},
    @{old = "13×20=260"; new = "90×76=6840"},
    @{old = "11×66=726"; new = "91×47=4277"},
    @{old = "57×12=684"; new = "26×84=2184"},
    @{old = "61×69=4209"; new = "31×70=2170"},
    @{old = "51×13=663"; new = "66×86=5676"},
    @{old = "19×30=570"; new = "97×20=1940"},
    @{old = "45×17=765"; new = "75×67=5025"},
    @{old = "38×78=2964"; new = "41×51=2091"},
    @{old = "59×66=3894"; new = "82×33=2706"},
    @{old = "21×71=1491"; new = "37×68=2516"},
    @{old = "21×87=1827"; new = "32×66=2112"},
    @{old = "11×18=198"; new = "86×99=8514"},
    @{old = "54×50=2700"; new = "48×24=1152"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
